# Apply "generic retrieval pipeline update" changes to Sheet1.
# Adds a new column E "Average annual surface temperature (tas)" with
# numeric-looking text values, updates a few existing cells (B3, D3, B4, D5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cells -------------------------------------------------

$ws.Range("B3").Value = "Wetlands (inland) - Permanent Rivers/Streams/Creeks (includes waterfalls), Wetlands (inland) - Permanent Freshwater Marshes/Pools (under 8ha), Tamaulipas, Texas, Coahuila, Nuevo León, Mexico, United States"
$ws.Range("D3").Value = "1.5 m"

$ws.Range("B4").Value = "Marine Neritic - Seagrass (Submerged), Marine Neritic - Subtidal Loose Rock/pebble/gravel, Marine Neritic - Subtidal Rock and Rocky Reefs, Norway, Portugal, Spain, Azores, Madeira, Canary Islands, France, Germany, Italy, Greece, Turkey, Bulgaria, Lebanon, Tunisia, Morocco, Malta, Jersey, Guernsey, Gibraltar, Denmark, Croatia, Belgium, Albania, Monaco, Montenegro, Netherlands, Romania, Slovenia, Algeria, Egypt, Libya"

$ws.Range("D5").Value = "7.375 kg"

# --- New column E header ----------------------------------------------------
# Copy formatting (bold font, border, centered alignment) from the D1 header
# cell, then overwrite with the new header text.

$ws.Range("D1").Copy($ws.Range("E1"))
$ws.Range("E1").Value = "Average annual surface temperature (tas)"

# --- New column E data cells -------------------------------------------------
# Values look numeric ("10.36" etc.) - force a text number format before
# assignment so they are stored as text (matching inlineStr/string type in
# the source data), then reset the style back to the default (no explicit
# style index), matching the rest of the data column cells.

$dataCells = @("E2", "E3", "E4", "E5", "E6", "E7", "E8")
foreach ($addr in $dataCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("E2").Value = "10.36"
$ws.Range("E3").Value = "21.78"
$ws.Range("E4").Value = "15.01"
$ws.Range("E5").Value = "22.48"
$ws.Range("E6").Value = "12.23"
$ws.Range("E7").Value = "-"
$ws.Range("E8").Value = "12.98"

$ws.Range("E2:E8").Style = "Normal"

# --- Dimension/used range is recalculated automatically by the engine, but
# touch it explicitly to be safe.
$ws.Range("A1:E8").Select()
